$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text so numeric-looking strings
# (e.g. "2.45", "230.83") are stored as text, matching the source data,
# not auto-converted to numbers by Excel type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.786.39'
$ws.Range("E2").Value = '  +4.86%  '

$ws.Range("D3").Value = '2.260.04'
$ws.Range("E3").Value = '  +2.45%  '

$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = '230.83'
$ws.Range("E5").Value = '  +0.45%  '

$ws.Range("E6").Value = '  +2.47%  '

$ws.Range("D7").Value = '63.08'
$ws.Range("E7").Value = '  +4.33%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").Value = '0.427'
$ws.Range("E9").Value = '  +6.53%  '

$ws.Range("D10").Value = '0.101'
$ws.Range("E10").Value = '  +13.57%  '

$ws.Range("D11").Value = '56.30'
$ws.Range("E11").Value = '  -1.17%  '

$ws.Range("D12").Value = '25.95'
$ws.Range("E12").Value = '  +17.42%  '

$ws.Range("E13").Value = '  +2.50%  '

$ws.Range("D14").Value = '2.598.68'
$ws.Range("E14").Value = '  +2.53%  '

$ws.Range("D15").Value = '15.67'
$ws.Range("E15").Value = '  +1.86%  '

$ws.Range("D16").Value = '5.92'
$ws.Range("E16").Value = '  +6.22%  '

$ws.Range("D17").Value = '0.822'
$ws.Range("E17").Value = '  +3.51%  '

$ws.Range("D18").Value = '2.289.44'
$ws.Range("E18").Value = '  +3.43%  '

$ws.Range("D19").Value = '43.724.55'
$ws.Range("E19").Value = '  +4.99%  '

$ws.Range("D20").Value = '0.0000102'
$ws.Range("E20").Value = '  +12.82%  '

$ws.Range("D21").Value = '73.66'
$ws.Range("E21").Value = '  +2.43%  '

$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("D23").Value = '255.91'
$ws.Range("E23").Value = '  +5.99%  '

$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("D25").Value = '2.45'
$ws.Range("E25").Value = '  +4.60%  '

$ws.Range("E26").Value = '  +1.92%  '

$ws.Range("D27").Value = '9.94'
$ws.Range("E27").Value = '  +3.62%  '

$ws.Range("D28").Value = '171.56'
$ws.Range("E28").Value = '  +1.88%  '

$ws.Range("D29").Value = '20.80'
$ws.Range("E29").Value = '  +5.57%  '

$ws.Range("E30").Value = '  -1.78%  '

$ws.Range("D31").Value = '2.84'
$ws.Range("E31").Value = '  +9.71%  '

$ws.Range("D32").Value = '1.39'
$ws.Range("E32").Value = '  -3.47%  '

$ws.Range("D33").Value = '0.123'
$ws.Range("E33").Value = '  +2.45%  '

$ws.Range("D34").Value = '0.0679'
$ws.Range("E34").Value = '  +5.59%  '

$ws.Range("E35").Value = '  +2.40%  '

$ws.Range("D36").Value = '4.92'
$ws.Range("E36").Value = '  -1.25%  '

$ws.Range("D37").Value = '3.86'
$ws.Range("E37").Value = '  +8.51%  '

$ws.Range("D38").Value = '6.70'
$ws.Range("E38").Value = '  +6.91%  '

$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("D40").Value = '0.0254'
$ws.Range("E40").Value = '  +5.57%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").Value = '8.37'
$ws.Range("E42").Value = '  -4.00%  '

$ws.Range("D43").Value = '17.43'
$ws.Range("E43").Value = '  +8.50%  '

$ws.Range("D44").Value = '0.0960'
$ws.Range("E44").Value = '  +0.63%  '

$ws.Range("D45").Value = '97.15'
$ws.Range("E45").Value = '  +0.91%  '

$ws.Range("E46").Value = '  -0.60%  '

$ws.Range("B47").Value = 'TerraClassic'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D47").Value = '0.000209'
$ws.Range("E47").Value = '  -14.66%  '

$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.455.15'
$ws.Range("E48").Value = '  -0.41%  '

$ws.Range("B49").Value = 'FTXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D49").Value = '4.31'
$ws.Range("E49").Value = '  -1.32%  '

$ws.Range("E50").Value = '  +5.10%  '

$ws.Range("E51").Value = '  +0.33%  '

# Restore the original (default/Normal) style on column D now that
# the text values are committed, so formatting matches the source.
$ws.Range("D2:D51").Style = "Normal"
